# reporte por sedes listo
# The "sede" (campus) labels in column A were re-typed from ALL CAPS to
# Title Case (e.g. "ALAMEDA" -> "Alameda"), and the view was scrolled /
# zoomed / re-selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sedes = [ordered]@{
    3  = "Alameda"
    4  = "Antonio Varas"
    5  = "Campus Arauco"
    6  = "Campus Villarrica"
    7  = "Concepción"
    8  = "Maipú"
    9  = "Melipilla"
    10 = "Nacimiento"
    11 = "Online"
    12 = "Padre Alonso De Ovalle"
    13 = "Plaza Norte"
    14 = "Plaza Oeste"
    15 = "Plaza Vespucio"
    16 = "Puente Alto"
    17 = "Puerto Montt"
    18 = "San Bernardo"
    19 = "San Carlos De Apoquindo"
    20 = "San Joaquín"
    21 = "Valparaíso"
    22 = "Viña Del Mar"
}

foreach ($row in $sedes.Keys) {
    $ws.Range("A$row").Value = $sedes[$row]
}

# Match the saved view state: zoomed to 131%, scrolled down, new selection.
$win = $excel.ActiveWindow
$win.Zoom = 131
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("L10").Select() | Out-Null
